# invitation task: used different path versions
# Add a new "path_version" column (J) with value 3 for every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 10).Value = "path_version"

$lastRow = 29
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 10).Value = 3
}

$ws.Range("K28").Select()
